$d = $word.ActiveDocument

# ------------------------------------------------------------------
# 1) "Per-class IoU values are **not computed** in this report." ->
#    "Per-class IoU from local test masks (`results/test_evaluation_metrics.txt`):"
# ------------------------------------------------------------------
$d.Content.Find.Execute(
    'Per-class IoU values are **not computed** in this report.',
    $true, $false, $false, $false, $false, $true, 1, $false,
    'Per-class IoU from local test masks (`results/test_evaluation_metrics.txt`):',
    2) | Out-Null

# ------------------------------------------------------------------
# 2) Locate the paragraph that used to explain how to generate
#    per-class metrics, and turn it + the paragraphs that follow it
#    into the new per-class IoU table and mean-IoU summary line.
# ------------------------------------------------------------------
$anchor = $d.Content
$anchor.Find.Execute(
    'Use `python dataset/test_segmentation.py --data_dir dataset/Offroad_Segmentation_testImages` to generate per-class metrics.'
) | Out-Null
$tableStartIndex = $anchor.Paragraphs.Item(1).Index

# First line of the new block replaces the found paragraph's text.
$d.Paragraphs.Item($tableStartIndex).Range.Text = '| Class | IoU |'

$newLines = @(
    '|---|---|',
    '| Background | N/A |',
    '| Trees | 0.0551 |',
    '| Lush Bushes | 0.0068 |',
    '| Dry Grass | 0.1271 |',
    '| Dry Bushes | 0.0524 |',
    '| Ground Clutter | 0.0318 |',
    '| Flowers | 0.0315 |',
    '| Logs | 0.0050 |',
    '| Rocks | 0.2125 |',
    '| Landscape | 0.5209 |',
    '| Sky | 0.9497 |',
    '',
    'Mean IoU (test set): **0.1996**'
)

$idx = $tableStartIndex
foreach ($line in $newLines) {
    $d.Paragraphs.Item($idx).Range.InsertParagraphAfter() | Out-Null
    $idx = $idx + 1
    $d.Paragraphs.Item($idx).Range.Text = $line
}

# ------------------------------------------------------------------
# 3) CPU benchmark number update: ~1204.7 -> ~3188.7
# ------------------------------------------------------------------
$d.Content.Find.Execute(
    '- **CPU:** ~1204.7 ms/image (fails <50 ms requirement on CPU, from `results/inference_benchmark.json`)',
    $true, $false, $false, $false, $false, $true, 1, $false,
    '- **CPU:** ~3188.7 ms/image (fails <50 ms requirement on CPU, from `results/inference_benchmark.json`)',
    2) | Out-Null

Write-Output "done"
